$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 20:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1418504
$ws.Range("C4").Value = 9868
$ws.Range("D4").Value = 300439
$ws.Range("E4").Value = 1033878
$ws.Range("G4").Value = 762
$ws.Range("H4").Value = 84187

# Row 15 - India
$ws.Range("D15").Value = 26392
$ws.Range("E15").Value = 49098

# Row 37 - Rumania
$ws.Range("E37").Value = 7005
$ws.Range("G37").Value = 34
$ws.Range("H37").Value = 1036

# Row 76 - Guinea
$ws.Range("B76").Value = 2374
$ws.Range("C76").Value = 76
$ws.Range("D76").Value = 856
$ws.Range("E76").Value = 1504
$ws.Range("G76").Value = 3
$ws.Range("H76").Value = 14

# Row 137 - Montenegro
$ws.Range("D137").Value = 307
$ws.Range("E137").Value = 8

# Row 140 - Ruanda
$ws.Range("B140").Value = 287
$ws.Range("C140").Value = 1
$ws.Range("D140").Value = 164
$ws.Range("E140").Value = 123
